# Natmi following Dr Hou advice:
# Recompute Wnt7b-Fzd1 LR-pair stats. The "Sending cluster" column now also
# includes "FAPs" (in addition to the original "sCs"), so the sheet grows
# from 5 data rows (sCs -> {ECs,FAPs,M1,M2,sCs}) to 10 data rows
# (FAPs -> {...} then sCs -> {...}), and every numeric column is refreshed
# with the new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A..D = Sending cluster, Ligand symbol, Receptor symbol, Target cluster
#          E..T = the 16 numeric NATMI metrics
$data = @(
    @("FAPs","Wnt7b","Fzd1","ECs",1,0.3333333333333333,0.072919,0.218757,0.1016383815134179,0.1016383815134179,3,1,1.269979666666667,3.809939,0.05093337778110388,0.05093337778110388,0.09260564731366666,0.833450825823,0.00517678608268288,0.005176786082682881),
    @("FAPs","Wnt7b","Fzd1","FAPs",1,0.3333333333333333,0.072919,0.218757,0.1016383815134179,0.1016383815134179,3,1,20.604156,61.812468,0.8263433572627787,0.8263433572627789,1.502434451364,13.521910062276,0.08398820140655291,0.08398820140655293),
    @("FAPs","Wnt7b","Fzd1","M1",1,0.3333333333333333,0.072919,0.218757,0.1016383815134179,0.1016383815134179,3,1,0.6519523333333334,1.955857,0.02614698121592406,0.02614698121592406,0.04753971219433333,0.427857409749,0.002657536852248261,0.002657536852248262),
    @("FAPs","Wnt7b","Fzd1","M2",1,0.3333333333333333,0.072919,0.218757,0.1016383815134179,0.1016383815134179,3,1,0.7481240000000001,2.244372,0.03000400976428538,0.03000400976428539,0.05455245395600001,0.4909720856040001,0.003049558991354755,0.003049558991354755),
    @("FAPs","Wnt7b","Fzd1","sCs",1,0.3333333333333333,0.072919,0.218757,0.1016383815134179,0.1016383815134179,3,1,1.659922,4.979766,0.06657227397590788,0.0665722739759079,0.121039852318,1.089358670862,0.006766298180579109,0.006766298180579111),
    @("sCs","Wnt7b","Fzd1","ECs",3,1,0.6445166666666666,1.93355,0.8983616184865821,0.898361618486582,3,1,1.269979666666667,3.809939,0.05093337778110388,0.05093337778110388,0.8185230614944444,7.36670755345,0.045756591698421,0.045756591698421),
    @("sCs","Wnt7b","Fzd1","FAPs",3,1,0.6445166666666666,1.93355,0.8983616184865821,0.898361618486582,3,1,20.604156,61.812468,0.8263433572627787,0.8263433572627789,13.2797219446,119.5174975014,0.7423551558562259,0.7423551558562259),
    @("sCs","Wnt7b","Fzd1","M1",3,1,0.6445166666666666,1.93355,0.8983616184865821,0.898361618486582,3,1,0.6519523333333334,1.955857,0.02614698121592406,0.02614698121592406,0.4201941447055556,3.78174730235,0.0234894443636758,0.0234894443636758),
    @("sCs","Wnt7b","Fzd1","M2",3,1,0.6445166666666666,1.93355,0.8983616184865821,0.898361618486582,3,1,0.7481240000000001,2.244372,0.03000400976428538,0.03000400976428539,0.4821783867333334,4.3396054806,0.02695445077293063,0.02695445077293063),
    @("sCs","Wnt7b","Fzd1","sCs",3,1,0.6445166666666666,1.93355,0.8983616184865821,0.898361618486582,3,1,1.659922,4.979766,0.06657227397590788,0.0665722739759079,1.069847394366666,9.628626549299998,0.05980597579532878,0.05980597579532879)

)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2          # data starts on row 2 (row 1 is the header)
    $row = $data[$i]
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
